# Apply updated TPM-derived values to LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5429463333333333
$ws.Range("H2").Value = 1.628839
$ws.Range("I2").Value = 0.04659251079363984
$ws.Range("J2").Value = 0.04659251079363985
$ws.Range("M2").Value = 49.89274333333334
$ws.Range("N2").Value = 149.67823
$ws.Range("O2").Value = 0.8663408689480834
$ws.Range("P2").Value = 0.8663408689480835
$ws.Range("Q2").Value = 27.08908205277444
$ws.Range("R2").Value = 243.80173847497
$ws.Range("S2").Value = 0.0403649962874349
$ws.Range("T2").Value = 0.04036499628743491

# Row 3
$ws.Range("G3").Value = 0.5429463333333333
$ws.Range("H3").Value = 1.628839
$ws.Range("I3").Value = 0.04659251079363984
$ws.Range("J3").Value = 0.04659251079363985
$ws.Range("O3").Value = 0.06984725491313053
$ws.Range("P3").Value = 0.06984725491313053
$ws.Range("Q3").Value = 2.184011036903111
$ws.Range("R3").Value = 19.656099332128
$ws.Range("S3").Value = 0.003254358978446148
$ws.Range("T3").Value = 0.003254358978446148

# Row 4
$ws.Range("G4").Value = 0.5429463333333333
$ws.Range("H4").Value = 1.628839
$ws.Range("I4").Value = 0.04659251079363984
$ws.Range("J4").Value = 0.04659251079363985
$ws.Range("M4").Value = 1.266267666666667
$ws.Range("N4").Value = 3.798803
$ws.Range("O4").Value = 0.02198755485004457
$ws.Range("P4").Value = 0.02198755485004457
$ws.Range("Q4").Value = 0.6875153866352223
$ws.Range("R4").Value = 6.187638479717
$ws.Range("S4").Value = 0.00102445538667645
$ws.Range("T4").Value = 0.00102445538667645

# Row 5
$ws.Range("G5").Value = 0.5429463333333333
$ws.Range("H5").Value = 1.628839
$ws.Range("I5").Value = 0.04659251079363984
$ws.Range("J5").Value = 0.04659251079363985
$ws.Range("M5").Value = 0.2206823333333333
$ws.Range("N5").Value = 0.6620469999999999
$ws.Range("O5").Value = 0.003831942516052412
$ws.Range("P5").Value = 0.003831942516052413
$ws.Range("Q5").Value = 0.1198186637147778
$ws.Range("R5").Value = 1.078367973433
$ws.Range("S5").Value = 0.0001785398230397794
$ws.Range("T5").Value = 0.0001785398230397795

# Row 6
$ws.Range("G6").Value = 0.5429463333333333
$ws.Range("H6").Value = 1.628839
$ws.Range("I6").Value = 0.04659251079363984
$ws.Range("J6").Value = 0.04659251079363985
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.187988666666667
$ws.Range("N6").Value = 6.563966
$ws.Range("O6").Value = 0.03799237877268909
$ws.Range("P6").Value = 0.03799237877268909
$ws.Range("Q6").Value = 1.187960423941556
$ws.Range("R6").Value = 10.691643815474
$ws.Range("S6").Value = 0.001770160318042569
$ws.Range("T6").Value = 0.00177016031804257

# Row 7
$ws.Range("I7").Value = 0.8858267105024722
$ws.Range("J7").Value = 0.8858267105024723
$ws.Range("M7").Value = 49.89274333333334
$ws.Range("N7").Value = 149.67823
$ws.Range("O7").Value = 0.8663408689480834
$ws.Range("P7").Value = 0.8663408689480835
$ws.Range("Q7").Value = 515.0233811528434
$ws.Range("R7").Value = 4635.21043037559
$ws.Range("S7").Value = 0.767427882114134
$ws.Range("T7").Value = 0.7674278821141343

# Row 8
$ws.Range("I8").Value = 0.8858267105024722
$ws.Range("J8").Value = 0.8858267105024723
$ws.Range("O8").Value = 0.06984725491313053
$ws.Range("P8").Value = 0.06984725491313053
$ws.Range("S8").Value = 0.06187256405732605
$ws.Range("T8").Value = 0.06187256405732606

# Row 9
$ws.Range("I9").Value = 0.8858267105024722
$ws.Range("J9").Value = 0.8858267105024723
$ws.Range("M9").Value = 1.266267666666667
$ws.Range("N9").Value = 3.798803
$ws.Range("O9").Value = 0.02198755485004457
$ws.Range("P9").Value = 0.02198755485004457
$ws.Range("Q9").Value = 13.07118854487767
$ws.Range("R9").Value = 117.640696903899
$ws.Range("S9").Value = 0.01947716338480766
$ws.Range("T9").Value = 0.01947716338480766

# Row 10
$ws.Range("I10").Value = 0.8858267105024722
$ws.Range("J10").Value = 0.8858267105024723
$ws.Range("M10").Value = 0.2206823333333333
$ws.Range("N10").Value = 0.6620469999999999
$ws.Range("O10").Value = 0.003831942516052412
$ws.Range("P10").Value = 0.003831942516052413
$ws.Range("Q10").Value = 2.278017881572333
$ws.Range("R10").Value = 20.502160934151
$ws.Range("S10").Value = 0.003394437033829275
$ws.Range("T10").Value = 0.003394437033829276

# Row 11
$ws.Range("I11").Value = 0.8858267105024722
$ws.Range("J11").Value = 0.8858267105024723
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.187988666666667
$ws.Range("N11").Value = 6.563966
$ws.Range("O11").Value = 0.03799237877268909
$ws.Range("P11").Value = 0.03799237877268909
$ws.Range("Q11").Value = 22.58575587840867
$ws.Range("R11").Value = 203.271802905678
$ws.Range("S11").Value = 0.03365466391237513
$ws.Range("T11").Value = 0.03365466391237513

# Row 12
$ws.Range("G12").Value = 0.7875243333333334
$ws.Range("H12").Value = 2.362573
$ws.Range("I12").Value = 0.06758077870388791
$ws.Range("J12").Value = 0.06758077870388793
$ws.Range("M12").Value = 49.89274333333334
$ws.Range("N12").Value = 149.67823
$ws.Range("O12").Value = 0.8663408689480834
$ws.Range("P12").Value = 0.8663408689480835
$ws.Range("Q12").Value = 39.29174943175445
$ws.Range("R12").Value = 353.6257448857901
$ws.Range("S12").Value = 0.05854799054651438
$ws.Range("T12").Value = 0.0585479905465144

# Row 13
$ws.Range("G13").Value = 0.7875243333333334
$ws.Range("H13").Value = 2.362573
$ws.Range("I13").Value = 0.06758077870388791
$ws.Range("J13").Value = 0.06758077870388793
$ws.Range("O13").Value = 0.06984725491313053
$ws.Range("P13").Value = 0.06984725491313053
$ws.Range("Q13").Value = 3.167830281255111
$ws.Range("R13").Value = 28.510472531296
$ws.Range("S13").Value = 0.004720331877358322
$ws.Range("T13").Value = 0.004720331877358323

# Row 14
$ws.Range("G14").Value = 0.7875243333333334
$ws.Range("H14").Value = 2.362573
$ws.Range("I14").Value = 0.06758077870388791
$ws.Range("J14").Value = 0.06758077870388793
$ws.Range("M14").Value = 1.266267666666667
$ws.Range("N14").Value = 3.798803
$ws.Range("O14").Value = 0.02198755485004457
$ws.Range("P14").Value = 0.02198755485004457
$ws.Range("Q14").Value = 0.9972166000132224
$ws.Range("R14").Value = 8.974949400119002
$ws.Range("S14").Value = 0.001485936078560459
$ws.Range("T14").Value = 0.00148593607856046

# Row 15
$ws.Range("G15").Value = 0.7875243333333334
$ws.Range("H15").Value = 2.362573
$ws.Range("I15").Value = 0.06758077870388791
$ws.Range("J15").Value = 0.06758077870388793
$ws.Range("M15").Value = 0.2206823333333333
$ws.Range("N15").Value = 0.6620469999999999
$ws.Range("O15").Value = 0.003831942516052412
$ws.Range("P15").Value = 0.003831942516052413
$ws.Range("Q15").Value = 0.1737927074367778
$ws.Range("R15").Value = 1.564134366931
$ws.Range("S15").Value = 0.0002589656591833575
$ws.Range("T15").Value = 0.0002589656591833576

# Row 16
$ws.Range("G16").Value = 0.7875243333333334
$ws.Range("H16").Value = 2.362573
$ws.Range("I16").Value = 0.06758077870388791
$ws.Range("J16").Value = 0.06758077870388793
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.187988666666667
$ws.Range("N16").Value = 6.563966
$ws.Range("O16").Value = 0.03799237877268909
$ws.Range("P16").Value = 0.03799237877268909
$ws.Range("Q16").Value = 1.723094316057556
$ws.Range("R16").Value = 10.691643815474
$ws.Range("S16").Value = 0.00256755454227139
$ws.Range("T16").Value = 0.00256755454227139
